$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add slides link for session 03 (week 3, row 4, column E)
$ws.Range("E4").Value = "slides/slides.html#/session-03-autonomy-as-a-basic-psychological-need-the-framework-of-self-determination-theory"

# Add prep material for week 4 (row 5, column D)
$ws.Range("D5").Value = "prep/p04.html"

# Update the active selection to D6
$ws.Range("D6").Select()
